# Correct Start Time / End Time values for all data rows.
# Start Time: 10/01/2024 -> 11/01/2024
# End Time:   11/01/2025 -> 12/01/2025
#
# Values are stored as plain text (not real dates), so we must force
# Excel to keep them as literal text via a leading apostrophe and then
# reset the cell style back to Normal so no stray number-format/style
# gets attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $startCell = $ws.Cells.Item($r, 1)
    $endCell = $ws.Cells.Item($r, 2)

    if ([string]$startCell.Text -eq "10/01/2024") {
        $startCell.Value = "'11/01/2024"
        $startCell.Style = "Normal"
    }
    if ([string]$endCell.Text -eq "11/01/2025") {
        $endCell.Value = "'12/01/2025"
        $endCell.Style = "Normal"
    }
}
